# Auto-generated edit script: update cryptos list values (2023-04-26 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.404.13'
$ws.Range("E2").Value = '  +3.33%  '

$ws.Range("D3").Value = '1.865.61'
$ws.Range("E3").Value = '  +1.89%  '

$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").Value = '''337.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.86%  '

$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("D7").Value = '''0.4707'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.61%  '

$ws.Range("D8").Value = '''0.3968'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.58%  '

$ws.Range("D9").Value = '''47.61'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.48%  '

$ws.Range("D10").Value = '''0.08022'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.24%  '

$ws.Range("D11").Value = '''0.9972'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.77%  '

$ws.Range("E12").Value = '  +4.11%  '

$ws.Range("D13").Value = '''6.027'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.52%  '

$ws.Range("D14").Value = '1.859.14'
$ws.Range("E14").Value = '  -0.53%  '

$ws.Range("D15").Value = '''7.247'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.65%  '

$ws.Range("D16").Value = '''90.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.52%  '

$ws.Range("D17").Value = '''1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '''0.00001041'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.90%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '''0.06618'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.63%  '

$ws.Range("D20").Value = '''17.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.64%  '

$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").Value = '28.421.94'
$ws.Range("E22").Value = '  +3.51%  '

$ws.Range("D23").Value = '''5.464'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.30%  '

$ws.Range("D24").Value = '''11.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.28%  '

$ws.Range("D25").Value = '''2.269'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.60%  '

$ws.Range("D26").Value = '2.080.53'
$ws.Range("E26").Value = '  +0.28%  '

$ws.Range("D27").Value = '''160.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.41%  '

$ws.Range("D28").Value = '''19.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.68%  '

$ws.Range("D29").Value = '''2.112'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.34%  '

$ws.Range("D30").Value = '''5.472'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.32%  '

$ws.Range("D31").Value = '''119.53'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.05%  '

$ws.Range("D32").Value = '''0.9628'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.28%  '

$ws.Range("D33").Value = '''0.09518'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.51%  '

$ws.Range("D34").Value = '''3.596'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.75%  '

$ws.Range("D35").Value = '''1.377'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.47%  '

$ws.Range("D36").Value = '''5.353'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.15%  '

$ws.Range("D37").Value = '''0.06103'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.80%  '

$ws.Range("D38").Value = '''0.02246'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.22%  '

$ws.Range("D39").Value = '''8.282'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.22%  '

$ws.Range("D40").Value = '''1.178'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.90%  '

$ws.Range("D41").Value = '''0.5923'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.30%  '

$ws.Range("D42").Value = '''1.001'
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").Value = '''0.1875'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.86%  '

$ws.Range("D44").Value = '''10.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.66%  '

$ws.Range("D45").Value = '''1.260'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.22%  '

$ws.Range("D46").Value = '''0.5552'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.15%  '

$ws.Range("E47").Value = '  +1.58%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.07327'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.21%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '''1.950'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.41%  '

$ws.Range("D50").Value = '''2.054'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.55%  '

$ws.Range("D51").Value = '''111.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.64%  '

